$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-02-17"

# Update the header label for the 2022 YTD column
$ws.Range("I1").Value = "2022 (through 02-17)"

# Update data values per the diff
$ws.Range("I2").Value = 160
$ws.Range("I3").Value = 78
$ws.Range("H12").Value = 203
$ws.Range("H14").Value = 1852
